$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.874.07'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '3.406.07'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.18'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.64'
$ws.Range("E6").Value = '  -1.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  +6.21%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.736'
$ws.Range("E9").Value = '  +5.38%  '
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.61'
$ws.Range("E11").Value = '  +2.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000223'
$ws.Range("E12").Value = '  +33.27%  '
$ws.Range("E13").Value = '  +8.67%  '
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.38'
$ws.Range("E15").Value = '  +7.04%  '
$ws.Range("D16").Value = '3.951.47'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '3.416.63'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.56'
$ws.Range("E18").Value = '  +8.32%  '
$ws.Range("E19").Value = '  +6.32%  '
$ws.Range("D20").Value = '61.894.51'
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '449.33'
$ws.Range("E21").Value = '  +42.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.82'
$ws.Range("E22").Value = '  +8.28%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.21'
$ws.Range("E24").Value = '  +2.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.29'
$ws.Range("E25").Value = '  +3.37%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.21'
$ws.Range("E26").Value = '  +10.36%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.32'
$ws.Range("E27").Value = '  +13.97%  '
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.63'
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.05'
$ws.Range("E31").Value = '  +4.81%  '
$ws.Range("E32").Value = '  -2.59%  '
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.47'
$ws.Range("E34").Value = '  -4.13%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +3.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.89'
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.136'
$ws.Range("E39").Value = '  +8.10%  '
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.38'
$ws.Range("E43").Value = '  +11.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '143.85'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.59'
$ws.Range("E45").Value = '  +15.84%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("E47").Value = '  -2.01%  '
$ws.Range("E48").Value = '  +22.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.59'
$ws.Range("E49").Value = '  +5.80%  '
$ws.Range("E50").Value = '  +6.03%  '
$ws.Range("D51").Value = '3.750.51'
$ws.Range("E51").Value = '  -0.73%  '
